$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Data edits: rows 24-26 ("control") stimDir now point at the emotiv/robot
# clip instead of the old control_N.mp4 files (the three control mp4 shared
# strings become unused and are dropped by Excel on save).
$ws.Range("F24").Value = "./stimuli/robot_right_3.mp4"
$ws.Range("F25").Value = "./stimuli/robot_right_3.mp4"
$ws.Range("F26").Value = "./stimuli/robot_right_3.mp4"

# View state: scroll/zoom in on the edited rows and leave the selection on
# the last-touched cell.
$win = $excel.ActiveWindow
$win.Zoom = 130
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws.Range("F14").Select() | Out-Null
